$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("205:205").Insert()
$ws.Range("A205").Value = 5
$ws.Range("B205").Value = "Macroferia Regional de Talca"
$ws.Range("C205").Value = "Maule"
$ws.Range("D205").Value = 45097
$ws.Range("E205").Value = 7
$ws.Range("F205").Value = 100112017
$ws.Range("G205").Value = "Apio"
$ws.Range("H205").Value = "Americana (o)"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 700
$ws.Range("K205").Value = 5000
$ws.Range("L205").Value = 5000
$ws.Range("M205").Value = 5000
$ws.Range("N205").Value = "$/docena de matas"
$ws.Range("O205").Value = "Provincia del Elquí"
$ws.Range("P205").Value = 833
$ws.Range("Q205").Value = 6
$ws.Range("R205").Value = "Hortaliza"
